$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Fixing Legs For 1,2, and 3" -- correct the row-9 (Leg/R3) calibration
# readings: Combo B's "beta 0" (D9) and "alpha 45" (G9) measurements.
# H9 ("alpha 0") is a shared formula (=D9-2*(D9-G9)) and recalculates
# automatically once the inputs change.
$ws.Range("D9").Value = 1300
$ws.Range("G9").Value = 850

# Restore the saved cursor position / selection on the sheet.
$ws.Range("G11").Select()
